# Add a new "TGZpMotion" worksheet, cloned from "TGMcontroller", trimmed
# down to the subset of rows that apply to the extra ("TGMotion") structure
# that does not fall under the TGM tab, then tidy up TGMcontroller's sheet
# view / column formatting so it matches its sibling sheets.

$wb = $excel.ActiveWorkbook

# TGMcontroller is the template for the new sheet - clone it and place the
# copy right after the original, at the end of the tab strip.
$src = $wb.Worksheets.Item("TGMcontroller")
$src.Copy($null, $src)

$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "TGZpMotion"

# The new sheet only needs 17 rows (vs. 20 on TGMcontroller): drop the
# "Napájení" connector row (old row 16) and the feedback/encoder connector
# rows at the bottom (old rows 19-20), leaving CAN + I/O connector rows in
# place of them.
$new.Rows.Item(16).Delete()
$new.Rows.Item(18).Delete()
$new.Rows.Item(18).Delete()

# Match the selection left behind in the saved file. Copying the sheet
# already makes it the active tab (and drops tabSelected from the source
# sheet), so this is the only manual selection tweak needed.
$new.Range("B21").Select() | Out-Null
